$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2031.6129
$ws.Range("I70").Value = 3580
$ws.Range("J70").Value = 1493.0435
$ws.Range("K70").Value = 10740
$ws.Range("L70").Value = 4479.1305
$ws.Range("M70").Value = -10470
$ws.Range("N70").Value = -5019.1305

$ws.Range("H73").Value = 2031.6129
$ws.Range("I73").Value = 3580
$ws.Range("J73").Value = 1493.0435
$ws.Range("K73").Value = 10740
$ws.Range("L73").Value = 4479.1305
$ws.Range("M73").Value = -9804
$ws.Range("N73").Value = -6351.1305

$ws.Range("H129").Value = 3079.5
$ws.Range("I129").Value = 888.5
$ws.Range("J129").Value = 3517.7
$ws.Range("K129").Value = 2665.5
$ws.Range("L129").Value = 10553.1
$ws.Range("M129").Value = 2334.5
$ws.Range("N129").Value = -20553.1

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5720
$ws.Range("I28").Value = 5720
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5720
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -5528
$ws.Range("N28").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H63").Value = 2211.5
$ws.Range("I63").Value = 2211.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2211.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1525.5
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2211.5
$ws.Range("I66").Value = 2211.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11057.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7625.5
$ws.Range("N66").ClearContents()

$ws.Range("H88").Value = 4425.143
$ws.Range("I88").Value = 5494
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 5494
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -5088
$ws.Range("N88").Value = -3812

$ws.Range("H91").Value = 4425.143
$ws.Range("I91").Value = 5494
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 5494
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -4090
$ws.Range("N91").Value = -5808

$ws.Range("H99").Value = 5720
$ws.Range("I99").Value = 5720
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5720
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2725
$ws.Range("N99").ClearContents()

$ws.Range("H109").Value = 48000
$ws.Range("J109").Value = 48000
$ws.Range("L109").Value = 48000
$ws.Range("N109").Value = -50774

$ws.Range("H112").Value = 34000
$ws.Range("J112").Value = 34000
$ws.Range("L112").Value = 34000
$ws.Range("N112").Value = -36954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4029.2727
$ws.Range("J86").Value = 3992.6428
$ws.Range("L86").Value = 3992.6428
$ws.Range("N86").Value = -6238.6428

$ws.Range("H89").Value = 4029.2727
$ws.Range("J89").Value = 3992.6428
$ws.Range("L89").Value = 19963.214
$ws.Range("N89").Value = -31195.214

$ws.Range("H110").Value = 32851
$ws.Range("J110").Value = 32851
$ws.Range("L110").Value = 32851
$ws.Range("N110").Value = -41031

$ws.Range("H134").Value = 1981.25
$ws.Range("I134").Value = 1512.3334
$ws.Range("K134").Value = 4537.0002
$ws.Range("M134").Value = -2002.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26347966
$ws.Range("I31").Value = 100001730
$ws.Range("J31").Value = 43049.465
$ws.Range("K31").Value = 100001730
$ws.Range("L31").Value = 43049.465
$ws.Range("M31").Value = -100001435
$ws.Range("N31").Value = -43639.465

$ws.Range("H34").Value = 26347966
$ws.Range("I34").Value = 100001730
$ws.Range("J34").Value = 43049.465
$ws.Range("K34").Value = 100001730
$ws.Range("L34").Value = 43049.465
$ws.Range("M34").Value = -100001528
$ws.Range("N34").Value = -43453.465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 9459.208000000001
$ws.Range("I121").Value = 336
$ws.Range("J121").Value = 11860.053
$ws.Range("K121").Value = 1008
$ws.Range("L121").Value = 35580.159
$ws.Range("M121").Value = 302
$ws.Range("N121").Value = -38200.159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4067.5
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4601.25
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 4601.25
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -6597.25

$ws.Range("H83").Value = 4067.5
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4601.25
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 23006.25
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -32990.25

$ws.Range("H111").Value = 26633.334
$ws.Range("J111").Value = 26633.334
$ws.Range("L111").Value = 26633.334
$ws.Range("N111").Value = -32767.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1607
$ws.Range("I61").Value = 1150
$ws.Range("J61").Value = 2292.5
$ws.Range("K61").Value = 1150
$ws.Range("L61").Value = 2292.5
$ws.Range("M61").Value = -948
$ws.Range("N61").Value = -2696.5

$ws.Range("H110").Value = 39000
$ws.Range("J110").Value = 39000
$ws.Range("L110").Value = 39000
$ws.Range("N110").Value = -47180

$ws.Range("H113").Value = 1607
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 2292.5
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 2292.5
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -6632.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 19235.572
$ws.Range("J86").Value = 19235.572
$ws.Range("L86").Value = 19235.572
$ws.Range("N86").Value = -21481.572

$ws.Range("H89").Value = 19235.572
$ws.Range("J89").Value = 19235.572
$ws.Range("L89").Value = 96177.86
$ws.Range("N89").Value = -107409.86
